$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: new capacitor footprint association
$ws.Range("A5").Value = 402
$ws.Range("B5").Value = "C81"
$ws.Range("C5").Value = "C12"
$ws.Range("D5").Value = 603

# Row 6: new capacitor footprint association
$ws.Range("A6").Value = 402
$ws.Range("C6").Value = "C14"
$ws.Range("B6").Value = "C91"
$ws.Range("D6").Value = 603

# Row 7: LED footprint association (B/C only)
$ws.Range("B7").Value = "C231"
$ws.Range("C7").Value = "C15"

# Row 8: LED footprint association (B/C only)
$ws.Range("B8").Value = "C11"
$ws.Range("C8").Value = "C17"

# Row 9: LED footprint association (B/C only)
$ws.Range("B9").Value = "C21"
$ws.Range("C9").Value = "C18"

# Row 10: LED footprint association (B/C only)
$ws.Range("B10").Value = "C12"
$ws.Range("C10").Value = "C19"

# Row 15: a new BOM item added further down
$ws.Range("B15").Value = "C111"

# Restore the selection to match the saved worksheet view
$ws.Range("B15").Select()
